# This script reproduces a refresh of the cryptocurrency price table:
# most rows get updated Price (D) / Volume(1h) (E) values, and the
# Filecoin / ImmutableX rows (45-46) swap places.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value that must remain plain text (matches the source
# workbook, where every B/C/D/E data cell is an inline string) without Excel
# auto-converting numeric-looking text (e.g. "24.10", "0.0000191") into a
# Double, which would silently change its value/representation. We briefly
# switch the cell to Text format, assign, then restore the original style so
# no visible formatting/style residue is left behind.
function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

# Row 2
Set-TextValue "D2" "94.974.15"
Set-TextValue "E2" "  -1.30%  "

# Row 3
Set-TextValue "D3" "3.449.31"
Set-TextValue "E3" "  +3.93%  "

# Row 4
Set-TextValue "E4" "  +0.17%  "

# Row 5
Set-TextValue "D5" "239.11"
Set-TextValue "E5" "  -3.70%  "

# Row 6
Set-TextValue "D6" "640.82"
Set-TextValue "E6" "  -1.45%  "

# Row 7
Set-TextValue "E7" "  +6.34%  "

# Row 8
Set-TextValue "D8" "0.402"
Set-TextValue "E8" "  -3.71%  "

# Row 9
Set-TextValue "E9" "  +0.14%  "

# Row 10
Set-TextValue "E10" "  +1.87%  "

# Row 11
Set-TextValue "D11" "3.452.01"
Set-TextValue "E11" "  +4.10%  "

# Row 12
Set-TextValue "E12" "  -3.64%  "

# Row 13
Set-TextValue "D13" "41.57"
Set-TextValue "E13" "  +4.04%  "

# Row 14
Set-TextValue "E14" "  +0.75%  "

# Row 15
Set-TextValue "D15" "94.798.22"
Set-TextValue "E15" "  -1.24%  "

# Row 16
Set-TextValue "D16" "4.100.32"
Set-TextValue "E16" "  +4.17%  "

# Row 17
Set-TextValue "E17" "  +2.67%  "

# Row 18
Set-TextValue "D18" "8.44"
Set-TextValue "E18" "  -0.40%  "

# Row 19
Set-TextValue "D19" "3.463.21"
Set-TextValue "E19" "  +4.27%  "

# Row 20
Set-TextValue "D20" "17.79"
Set-TextValue "E20" "  +4.80%  "

# Row 21
Set-TextValue "D21" "11.37"
Set-TextValue "E21" "  +8.99%  "

# Row 22
Set-TextValue "D22" "0.504"
Set-TextValue "E22" "  -5.56%  "

# Row 23
Set-TextValue "D23" "501.92"
Set-TextValue "E23" "  +0.09%  "

# Row 24
Set-TextValue "E24" "  -5.24%  "

# Row 25
Set-TextValue "D25" "0.0000191"
Set-TextValue "E25" "  -2.38%  "

# Row 26
Set-TextValue "D26" "6.59"
Set-TextValue "E26" "  +0.74%  "

# Row 27
Set-TextValue "D27" "94.52"
Set-TextValue "E27" "  -1.15%  "

# Row 28
Set-TextValue "D28" "3.638.30"
Set-TextValue "E28" "  +4.05%  "

# Row 29
Set-TextValue "D29" "11.98"
Set-TextValue "E29" "  +0.19%  "

# Row 30
Set-TextValue "D30" "11.71"
Set-TextValue "E30" "  +6.89%  "

# Row 31
Set-TextValue "E31" "  +0.11%  "

# Row 32
Set-TextValue "D32" "2.74"
Set-TextValue "E32" "  +11.51%  "

# Row 33
Set-TextValue "E33" "  -3.45%  "

# Row 34
Set-TextValue "E34" "  -1.08%  "

# Row 35
Set-TextValue "D35" "30.84"
Set-TextValue "E35" "  +10.72%  "

# Row 36
Set-TextValue "D36" "0.998"
Set-TextValue "E36" "  -0.21%  "

# Row 37
Set-TextValue "E37" "  +3.89%  "

# Row 38
Set-TextValue "D38" "7.67"
Set-TextValue "E38" "  +1.84%  "

# Row 39
Set-TextValue "E39" "  -1.47%  "

# Row 40
Set-TextValue "D40" "521.28"
Set-TextValue "E40" "  +3.64%  "

# Row 41
Set-TextValue "E41" "  -0.03%  "

# Row 42
Set-TextValue "D42" "0.151"
Set-TextValue "E42" "  +0.38%  "

# Row 43
Set-TextValue "D43" "0.914"
Set-TextValue "E43" "  +10.94%  "

# Row 44
Set-TextValue "D44" "24.10"
Set-TextValue "E44" "  -0.98%  "

# Row 45
Set-TextValue "B45" "Filecoin"
Set-TextValue "C45" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D45" "5.65"
Set-TextValue "E45" "  +3.48%  "

# Row 46
Set-TextValue "B46" "ImmutableX"
Set-TextValue "C46" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D46" "1.69"
Set-TextValue "E46" "  +2.68%  "

# Row 47
Set-TextValue "D47" "0.0414"
Set-TextValue "E47" "  -4.04%  "

# Row 48
Set-TextValue "D48" "3.48"
Set-TextValue "E48" "  -4.24%  "

# Row 49
Set-TextValue "E49" "  +9.13%  "

# Row 50
Set-TextValue "D50" "53.46"
Set-TextValue "E50" "  +1.01%  "

# Row 51
Set-TextValue "E51" "  +2.30%  "
